# Add a new "Active Status" column (I) to the User Data sheet and mark
# the existing record as active (TRUE), matching the new layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "Active Status"
$ws.Range("I2").Value = $true

# Leave the selection where Excel would land after typing the new value.
$ws.Range("I3").Select()
